$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 with new data
$data = @(
    @("Food",       200,  "2025-06-14", "Pizza dinner with friends"),
    @("lunch",      120,  "2025-06-12", "N/A"),
    @("Groceries",  2500, "2025-06-10", "Monthly grocery shopping"),
    @("Food",       500,  "2025-06-10", "Lunch at restaurant"),
    @("Transport",  700,  "2025-06-08", "Bus and metro fare"),
    @("Ice cream",  50,   "2025-06-02", "N/A"),
    @("Rent",       1000, "2025-06-01", "June rent payment"),
    @("Electricity",400,  "2025-05-14", "Monthly electricity bill")
)

$row = 2
foreach ($entry in $data) {
    # Column C holds date-like text ("2025-06-14"); force text format so
    # Excel doesn't auto-convert it into a date serial number.
    $ws.Cells.Item($row, 3).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
